$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "NUMERO_Contrato"
$ws.Range("E11").Value = 'Número do Contrato registrado no Sistema Integrado de Administração Financeira de Minas Gerais (SIAFI-MG).  "NA" em caso de ausência de Contratos por serem de entrega imediata.'

$ws.Range("B12").Value = "URL_INTEGRA_Contrato"
$ws.Range("E12").Value = "Link para íntegra do termo do Contrato celebrado e seus eventuais termos aditivos ou modificativos."

$ws.Range("B13").Value = "CODIGO_ORGAO_Contrato"
$ws.Range("E13").Value = "Código do órgão ou entidade que formalizou o instrumento de contratação."

$ws.Range("B14").Value = "ORGAO_Contrato"
$ws.Range("E14").Value = "Nome do órgão ou entidade responsável pelo Contrato"

$ws.Range("E16").Value = "Data, no formato YYYY-MM-DD, de início da vigência do Contrato"
$ws.Range("E17").Value = "Data, no formato YYYY-MM-DD, de fim da vigência do Contrato"
$ws.Range("E18").Value = "Data, no formato YYYY-MM-DD, atualizada do fim da vigência do Contrato"

$ws.Range("E19").Value = "Número de identificação do fornecedor. Pessoa Jurídica – CNPJ. Pessoa física - CPF. Formato: (eg. xx.xxx.xxx/yyyy-zz e xxx.xxx.xxx-zz)."

$ws.Range("E23").Value = "Código do item de material ou serviço adquirido pelo Estado."
$ws.Range("E24").Value = "Descrição dos itens de materiais ou serviços adquiridos pelo Estado."
$ws.Range("E25").Value = "Código da unidade orçamentária vinculada ao item material/serviço."
$ws.Range("E26").Value = "Nome da Unidade Orçamentária vinculada ao item material/serviço."
